# Apply the "Added Jump, StorePC, JCond instructions" edit described by the
# commit. Operates on Sheet1 of the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Rename the "JumpBI / JumpFI" instruction-format header to "JCOND"
#    and update the JAL row (now the StorePC format row) in the small
#    instruction format legend (I10:M14).
# ---------------------------------------------------------------------
$ws.Range("I10").Value = "JCOND"

$ws.Range("I11").Value = "STOREPC"
$ws.Range("J11").Value = "OpCode "
$ws.Range("K11").Value = "OpExt"
$ws.Range("L11").Value = "xxxx"
$ws.Range("M11").Value = "StorePC"

$ws.Range("I12").Value = "Load"
$ws.Range("J12").Value = "OpCode"
$ws.Range("K12").Value = "OpExt"
$ws.Range("L12").Value = "Address"
$ws.Range("M12").Value = "StoreTo"

$ws.Range("I13").Value = "Store"
$ws.Range("J13").Value = "OpCode"
$ws.Range("K13").Value = "OpExt"
$ws.Range("L13").Value = "Address"
$ws.Range("M13").Value = "StoreTo"
$ws.Rows.Item(13).RowHeight = 14.9

# New row describing the JUMP instruction format, copy the look of the
# row above (Store) and then overwrite its contents.
$ws.Range("I13:M13").Copy()
$ws.Range("I14:M14").PasteSpecial(-4122)
$ws.Range("I14").Value = "JUMP "
$ws.Range("J14").Value = "OpCode"
$ws.Range("K14").Value = "OpExt"
$ws.Range("L14").Value = "Address"
$ws.Range("M14").Value = "xxxx"
$ws.Rows.Item(14).RowHeight = 14.9

# ---------------------------------------------------------------------
# 2. Non-ALU opcode table (rows 38-43): opcode 1100 now has two
#    sub-instructions (JUMP and STOREPC, distinguished by opext), so the
#    remaining "Open" opcode slots shift down by one row.
# ---------------------------------------------------------------------

# Row 39 (opcode 1100) becomes the JUMP instruction entry.
$ws.Range("F39").Value = "0000"
$ws.Range("G39").Value = "JUMP"

# Row 40 becomes the StorePC instruction entry sharing opcode 1100 with
# row 39, so merge E39:E40 and give E39/E40 the merged-cell look used by
# the other merged opcode cells (top-aligned, like E3/E33).
$ws.Range("E3").Copy()
$ws.Range("E39:E40").PasteSpecial(-4122)
$ws.Range("E39").Value = "1100"
$ws.Range("E40").ClearContents()
$ws.Range("E39:E40").Merge()

$ws.Range("F40").Value = "0001"
$ws.Range("G40").Value = "STOREPC"

# The remaining "Open" rows shift down: old 1101/1110/1111 become
# 1101/1110/1111 -> now sit in rows 41/42/43.
$ws.Range("E41").Value = "1101"
$ws.Range("E42").Value = "1110"

# Brand new row 43, cloned from row 42's formatting.
$ws.Range("E42:G42").Copy()
$ws.Range("E43:G43").PasteSpecial(-4122)
$ws.Range("I42").Copy()
$ws.Range("I43").PasteSpecial(-4122)
$ws.Range("E43").Value = "1111"
$ws.Range("F43").Value = "XXXX"
$ws.Range("G43").Value = "Open"
$ws.Rows.Item(43).RowHeight = 14.1
